$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "SCHEME_MEASURES": rename indicator codes MQMS0x -> MQME00x
# ---------------------------------------------------------------------------
$wsMeasures = $wb.Worksheets.Item("SCHEME_MEASURES")
$wsMeasures.Range("A2").Value = "MQME001"
$wsMeasures.Range("A3").Value = "MQME002"
$wsMeasures.Range("A4").Value = "MQME003"
$wsMeasures.Range("A5").Value = "MQME004"
$wsMeasures.Range("A6").Value = "MQME005"

# ---------------------------------------------------------------------------
# Sheet "METADATA_ISSUES": rename indicator code MQME12 -> MQME014
# ---------------------------------------------------------------------------
$wsIssues = $wb.Worksheets.Item("METADATA_ISSUES")
$wsIssues.Range("A2").Value = "MQME014"

# ---------------------------------------------------------------------------
# Sheet "METADATA_MEASURES": drop the "Total number of columns" row (old row
# 2) and shift the remaining two rows up, renumbering the indicator codes.
# ---------------------------------------------------------------------------
$wsMetaMeasures = $wb.Worksheets.Item("METADATA_MEASURES")
$wsMetaMeasures.Rows.Item(2).Delete()

$wsMetaMeasures.Range("A2").Value = "MQME006"
$wsMetaMeasures.Range("B2").Value = "Total number of length-required columns"
$wsMetaMeasures.Range("C2").Value = 11

$wsMetaMeasures.Range("A3").Value = "MQME007"
$wsMetaMeasures.Range("B3").Value = "Total number of NUMBER columns"
$wsMetaMeasures.Range("C3").Value = 27

# ---------------------------------------------------------------------------
# Sheet "METADATA_METRICS": replace the 7 IQME* indicators with the new set
# of 11 MQID* indicators.
# Column C holds percentages stored as plain text (e.g. "100.00%"), so force
# a text number format first to stop them being parsed as numeric percents.
# ---------------------------------------------------------------------------
$wsMetrics = $wb.Worksheets.Item("METADATA_METRICS")
$wsMetrics.Range("C2:C12").NumberFormat = "@"

$wsMetrics.Range("A2").Value = "MQID001"
$wsMetrics.Range("B2").Value = "Table names in singular"
$wsMetrics.Range("C2").Value = "100.00%"

$wsMetrics.Range("A3").Value = "MQID002"
$wsMetrics.Range("B3").Value = "Table with recommended name length"
$wsMetrics.Range("C3").Value = "100.00%"

$wsMetrics.Range("A4").Value = "MQID003"
$wsMetrics.Range("B4").Value = "Columns with correct prefixes"
$wsMetrics.Range("C4").Value = "97.83%"

$wsMetrics.Range("A5").Value = "MQID004"
$wsMetrics.Range("B5").Value = "Columns with recommended name size"
$wsMetrics.Range("C5").Value = "100.00%"

$wsMetrics.Range("A6").Value = "MQID005"
$wsMetrics.Range("B6").Value = "Columns with comments"
$wsMetrics.Range("C6").Value = "100.00%"

$wsMetrics.Range("A7").Value = "MQID006"
$wsMetrics.Range("B7").Value = "Table with standard PK prefixes"
$wsMetrics.Range("C7").Value = "100.00%"

$wsMetrics.Range("A8").Value = "MQID007"
$wsMetrics.Range("B8").Value = "Table with standard FK prefixes"
$wsMetrics.Range("C8").Value = "100.00%"

$wsMetrics.Range("A9").Value = "MQID008"
$wsMetrics.Range("B9").Value = "Table with standard UK prefixes"
$wsMetrics.Range("C9").Value = "0.00%"

$wsMetrics.Range("A10").Value = "MQID009"
$wsMetrics.Range("B10").Value = "NUMBER columns with valid scale"
$wsMetrics.Range("C10").Value = "100.00%"

$wsMetrics.Range("A11").Value = "MQID010"
$wsMetrics.Range("B11").Value = "Columns with valid num_distinct"
$wsMetrics.Range("C11").Value = "100.00%"

$wsMetrics.Range("A12").Value = "MQID011"
$wsMetrics.Range("B12").Value = "Columns with valid num_nulls"
$wsMetrics.Range("C12").Value = "100.00%"
